$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "2 Ambientes - Oportunidad"
$ws.Range("C3").Value = "Ideal Inversión / AirBnb"

$ws.Range("A5").Value = "Caballito"
$ws.Range("B5").Value = "USD 89.000"
$ws.Range("C5").Value = "3 Ambientes luminoso"
$ws.Range("D5").Value = "https://www.zonaprop.com.ar"
